# Apply "a working version with launch analysis" edit.
#
# Summary of changes (per the OOXML diff):
#  - Launch sheet: add a small "launch mux" breakdown table (rows 15-22,
#    columns E-H) with shared strings, literal values, and formulas
#    (including a shared formula for G20:G21), plus new column widths.
#  - Active sheet switches from "Receive" to "Launch" (selection F17).
#  - bookViews / activeTab follow the newly active sheet automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Launch")

# --- New mini "launch" analysis table -------------------------------------

$ws.Range("F15").Value = "1000 rounds * 4 ranks"

$ws.Range("E16").Value = "Components"
$ws.Range("F16").Value = "total time"

$ws.Range("E17").Value = "dpu_switch_mux_for_rank"

$ws.Range("E19").Value = "ufi_select_all_even_disabled"
$ws.Range("F19").Value = 0.000853539
$ws.Range("G19").Formula = "=F19/4000"
$ws.Range("H19").Formula = "=G19/G22"

$ws.Range("E20").Value = "ufi_set_mram_mux"
$ws.Range("F20").Value = 0.0456674
$ws.Range("G20:G21").Formula = "=F20/4000"
$ws.Range("H20").Formula = "=G20/G22"

$ws.Range("E21").Value = "dpu_check_wavegen_mux_status_for_rank"
$ws.Range("F21").Value = 0.221554
$ws.Range("H21").Formula = "=G21/G22"

$ws.Range("F22").Formula = "=SUM(F19:F21)"
$ws.Range("G22").Formula = "=SUM(G19:G21)"

# --- Column widths for the new columns (C, E, F:Q) -------------------------

$ws.Columns.Item(3).ColumnWidth = 10
$ws.Columns.Item(5).ColumnWidth = 36.166666666666664
for ($col = 6; $col -le 17; $col++) {
    $ws.Columns.Item($col).ColumnWidth = 23.666666666666668
}

# --- Window / selection: Launch becomes the active sheet -------------------

$ws.Activate()
[void]$ws.Range("F17").Select()
